# Auto-generated edit script for Bethune-Cookman University Organizations workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target table data: row 1 is the new header, rows 2-21 are the reordered/
# filtered organization rows (columns A (Organization Name) and B (Categories)
# swapped vs. the original; several promotional rows removed; column M dropped).
$data = @(
    @('Organization Name', 'Categories', 'Org URL', 'Image URL', 'Description', 'Email', 'Phone', 'Website', 'LinkedIn', 'Instagram', 'Facebook', 'Twitter'),
    @('Student Organizations', 'General', 'https://www.cookman.edu/studentexperience/student-organizations.html', '', 'African Students Association', '', '', '', '', '', '', ''),
    @('Greek-lettered organizations', 'Academic', 'https://www.cookman.edu/studentexperience/student-organizations.html', '', 'Fraternities and sororities are about friendships, scholarship, community service, brotherhood, sisterhood and leadership. They’re about bettering the B-CU community and the city of Daytona Beach through service and social functions. They’re about traditions, learning valuable lessons, and establishing life-long friendships. They’re about shaping the identity of a continued Greek Life presence since 1948.', '', '', '', '', '', '', ''),
    @('Greek Life: Info for Current Students', 'Greek Life', 'https://www.cookman.edu/studentexperience/_files/greek-forms/bcuregform.pdf', '', 'Membership Intake Seminar', '', '', '', '', '', '', ''),
    @('STUDENT EXPERIENCE', 'General', 'https://www.cookman.edu/studentexperience/index.html', '', '', '', '', '', '', '', '', ''),
    @('WHO WE ARE', 'General', 'https://www.cookman.edu/about/index.html', '', '', '', '', '', '', '', '', ''),
    @('WILDCAT WEB', 'General', 'https://www.cookman.edu/studentexperience/student-organizations.html', '', '', '', '', '', '', '', '', ''),
    @('Financial Reports', 'General', 'https://www.cookman.edu/aid/financial-reports.html', '', '', '', '', '', '', '', '', ''),
    @('SACSCOC RISE QEP', 'General', 'https://www.cookman.edu/qep/index.html', '', '', '', '', '', '', '', '', ''),
    @('Marketing and Communications', 'General', 'https://www.cookman.edu/comms/index.html', '', '', '', '', '', '', '', '', ''),
    @('B-CU Jobs', 'General', 'https://www.cookman.edu/studentexperience/student-organizations.html', '', '', '', '', '', '', '', '', ''),
    @('Payment Center', 'General', 'https://www.cookman.edu/payment-center/index.html', '', '', '', '', '', '', '', '', ''),
    @('Clery Report', 'General', 'https://www.cookman.edu/campussafety/_files/2024-safety-_-fire-report-final.pdf', '', '', '', '', '', '', '', '', ''),
    @('Accessibility Statement', 'General', 'https://www.cookman.edu/cit/ws/wa/web-accessibility-statement.html', '', '', '', '', '', '', '', '', ''),
    @('Open Bids', 'General', 'https://www.cookman.edu/open-bids/index.html', '', '', '', '', '', '', '', '', ''),
    @('Residence Life', 'General', 'https://www.cookman.edu/studentexperience/residence-life.html', '', '', '', '', '', '', '', '', ''),
    @('Center for Civic Engagement', 'General', 'https://www.cookman.edu/studentexperience/cce-wsc.html', '', '', '', '', '', '', '', '', ''),
    @('Chaplaincy & Religious Life', 'Religious', 'https://www.cookman.edu/crl/index.html', '', '', '', '', '', '', '', '', ''),
    @('Future Students', 'General', 'https://www.cookman.edu/prospective/index.html', '', '', '', '', '', '', '', '', ''),
    @('Current Students', 'General', 'https://www.cookman.edu/currentstudents/index.html', '', '', '', '', '', '', '', '', ''),
    @('Student Organization Program / Event Approval', 'General', 'https://www.cookman.edu/studentexperience/student-organizations.html', '', 'B-CU students please submit your request for approval of your Programs / Events below.', '', '', '', '', '', '', ''),
)

$rowCount = $data.Count
$colCount = $data[0].Count

$values = New-Object "object[,]" $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $values[$r, $c] = $data[$r][$c]
    }
}

# Write the full A1:L21 block in one shot (keeps the existing header style on row 1 intact).
$target = $ws.Range("A1").Resize($rowCount, $colCount)
$target.Value = $values

# Wipe everything outside the new A1:L21 extent so the used range (and saved
# <dimension>) shrinks down from the old A1:M27 to A1:L21, matching the diff.
$ws.Range("M1:M27").Clear()
$ws.Range("A22:L27").Clear()

# Column widths: A/B swap to (47, 12); G:L narrow down now that M (Tiktok) is gone.
$ws.Columns.Item(1).ColumnWidth = 47
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 50
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 50
$ws.Columns.Item(6).ColumnWidth = 7
$ws.Columns.Item(7).ColumnWidth = 7
$ws.Columns.Item(8).ColumnWidth = 9
$ws.Columns.Item(9).ColumnWidth = 10
$ws.Columns.Item(10).ColumnWidth = 11
$ws.Columns.Item(11).ColumnWidth = 10
$ws.Columns.Item(12).ColumnWidth = 9

$ws.Range("A1").Select()

Write-Output $ws.UsedRange.Address()
